# Update the "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) columns
# on the Training Dashboard sheet for the new progress date of 04-Nov-2025:
#   - PERIOD TO EXPIRE decreases by 1 (one more day has elapsed)
#   - LAST UPDATE changes from 03-Nov-2025 to 04-Nov-2025
#
# Column I is stored as literal text (not a real date), so a helper cell
# formatted as Text ("@") is used with Copy / PasteSpecial(values) to push
# the new literal string in without Excel's automatic text->date
# conversion touching the target cell's number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$newLastUpdate = "04-Nov-2025"

# Scratch cell, far away from the used range, used only to stage the
# literal text value for a format-preserving paste into column I.
$helper = $ws.Cells.Item(1000, 1)
$helper.NumberFormat = "@"
$helper.Value = $newLastUpdate
$helper.Copy()

for ($row = 3; $row -le 17; $row++) {
    $periodCell = $ws.Cells.Item($row, 8)
    $periodCell.Value = $periodCell.Value2 - 1

    $lastUpdateCell = $ws.Cells.Item($row, 9)
    $lastUpdateCell.PasteSpecial(-4163)  # xlPasteValues
}

$helper.Clear()
$excel.CutCopyMode = $false
